$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These text values must be stored as TEXT (shared string), not auto-converted
# to numbers/dates by Excel's input parser. Using a leading apostrophe forces
# text entry; re-applying the "Normal" style afterwards drops the transient
# quote-prefix formatting so the cell keeps the sheet's default style.
$ws.Range("A2").Value = "ocds-twb234-0005"
$ws.Range("B2").Value = "'3568999"
$ws.Range("C2").Value = "Activo"
$ws.Range("D2").Value = "'2019-03-16"
$ws.Range("E2").Value = "'30628707093"
$ws.Range("F2").Value = "HAL2000"
$ws.Range("G2").Value = "ARS"
$ws.Range("H2").Value = "'2000000"

$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("H2").Style = "Normal"

# Align page margins with the target layout (values are in points; 72pt = 1in)
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
